$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Se cambio la descripción"
$tr.LanguageID = "es-CO"
